# The deck ships two DrawingML themes:
#   ppt/theme/theme1.xml -> bound to the slide master ("Integral" / Red Violet)
#   ppt/theme/theme2.xml -> bound to the notes master  ("Office Theme")
#
# The target edit swaps the two theme's colour schemes (font scheme and
# format scheme are identical between them, so only the 12 theme colours -
# and the scheme's display name - actually change). We drive this through
# the real PowerPoint automation surface: Master.Theme.ThemeColorScheme,
# which maps index 1-12 onto dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
# (the standard MSO theme-colour order) and writes straight back into the
# theme part that slide master's Design uses.

$p = $ppt.ActivePresentation
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

# Stock "Office" theme colours (what used to live in theme2.xml) replace
# the "Integral"/Red Violet colours current in theme1.xml.
$tcs.Colors(1).RGB  = 0          # dk1      000000
$tcs.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388    # dk2      44546A
$tcs.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407      # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308   # accent5  4472C4
$tcs.Colors(10).RGB = 4697456    # accent6  70AD47
$tcs.Colors(11).RGB = 12673797   # hlink    0563C1
$tcs.Colors(12).RGB = 7491477    # folHlink 954F72
